# "Se mejora el menu" - append new WhatsApp chat rows (Fecha / Mensaje)
# to the bottom of the "Chats" sheet, rows 31-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Fecha (column A), Mensaje (column B)
$rows = @(
    @(31, "03-05-2022 09:38", "Hola"),
    @(32, "03-05-2022 09:39", "Hola"),
    @(33, "03-05-2022 09:39", "hola"),
    @(34, "03-05-2022 09:41", "hola"),
    @(35, "03-05-2022 09:42", "hola"),
    @(36, "03-05-2022 09:42", "1"),
    @(37, "03-05-2022 09:48", "hola"),
    @(38, "03-05-2022 09:48", "hola"),
    @(39, "03-05-2022 09:48", "hola"),
    @(40, "03-05-2022 09:49", "1"),
    @(41, "03-05-2022 09:51", "hola"),
    @(42, "03-05-2022 09:52", "1"),
    @(43, "03-05-2022 09:53", "hola"),
    @(44, "03-05-2022 09:53", "1")
)

foreach ($row in $rows) {
    $rowNum = $row[0]
    $fecha = $row[1]
    $mensaje = $row[2]

    # Column A: timestamp text (never numeric-looking, stored as text naturally)
    $ws.Cells.Item($rowNum, 1).Value = $fecha

    # Column B: message text. Some messages are the bare digit "1", which
    # Excel would otherwise coerce to a number - force text storage for those.
    if ($mensaje -eq "1") {
        $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
        $ws.Cells.Item($rowNum, 2).Value = $mensaje
    } else {
        $ws.Cells.Item($rowNum, 2).Value = $mensaje
    }
}

# Drop the temporary text formatting again so the new cells stay styleless,
# matching the rest of the sheet.
$ws.Range("B31:B44").ClearFormats()
